$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking values
# (e.g. "1.019") are not auto-converted into floating point numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '27.938.83'
$ws.Range("E2").Value = '  +0.84%  '

# Row 3
$ws.Range("D3").Value = '1.882.05'
$ws.Range("E3").Value = '  +0.20%  '

# Row 4
$ws.Range("D4").Value = '1.019'
$ws.Range("E4").Value = '  +1.65%  '

# Row 5
$ws.Range("D5").Value = '335.42'
$ws.Range("E5").Value = '  +1.47%  '

# Row 6
$ws.Range("D6").Value = '1.018'
$ws.Range("E6").Value = '  +1.53%  '

# Row 7
$ws.Range("D7").Value = '0.4688'
$ws.Range("E7").Value = '  -0.69%  '

# Row 8
$ws.Range("D8").Value = '0.3912'
$ws.Range("E8").Value = '  -1.44%  '

# Row 9
$ws.Range("D9").Value = '47.53'
$ws.Range("E9").Value = '  -1.44%  '

# Row 10
$ws.Range("D10").Value = '0.07953'
$ws.Range("E10").Value = '  -1.15%  '

# Row 11
$ws.Range("D11").Value = '1.009'
$ws.Range("E11").Value = '  -1.62%  '

# Row 12
$ws.Range("D12").Value = '21.65'
$ws.Range("E12").Value = '  -0.99%  '

# Row 13
$ws.Range("D13").Value = '1.908.63'
$ws.Range("E13").Value = '  +3.05%  '

# Row 14
$ws.Range("D14").Value = '5.940'
$ws.Range("E14").Value = '  -0.52%  '

# Row 15
$ws.Range("D15").Value = '7.084'
$ws.Range("E15").Value = '  -1.28%  '

# Row 16
$ws.Range("D16").Value = '1.021'
$ws.Range("E16").Value = '  +1.63%  '

# Row 17
$ws.Range("D17").Value = '0.06776'
$ws.Range("E17").Value = '  +2.52%  '

# Row 18
$ws.Range("D18").Value = '87.06'
$ws.Range("E18").Value = '  -0.09%  '

# Row 19
$ws.Range("D19").Value = '0.00001044'
$ws.Range("E19").Value = '  -0.25%  '

# Row 20
$ws.Range("D20").Value = '17.04'
$ws.Range("E20").Value = '  -1.17%  '

# Row 21
$ws.Range("D21").Value = '1.017'
$ws.Range("E21").Value = '  +1.51%  '

# Row 22
$ws.Range("D22").Value = '27.945.39'
$ws.Range("E22").Value = '  +0.81%  '

# Row 23
$ws.Range("D23").Value = '5.472'
$ws.Range("E23").Value = '  -0.60%  '

# Row 24
$ws.Range("D24").Value = '10.93'
$ws.Range("E24").Value = '  -0.92%  '

# Row 25
$ws.Range("D25").Value = '2.365'
$ws.Range("E25").Value = '  +2.89%  '

# Row 26
$ws.Range("D26").Value = '2.126.28'
$ws.Range("E26").Value = '  +2.34%  '

# Row 27
$ws.Range("D27").Value = '160.03'
$ws.Range("E27").Value = '  +2.16%  '

# Row 28
$ws.Range("E28").Value = '  -1.89%  '

# Row 29
$ws.Range("D29").Value = '2.071'
$ws.Range("E29").Value = '  -1.20%  '

# Row 30
$ws.Range("D30").Value = '5.457'
$ws.Range("E30").Value = '  -2.37%  '

# Row 31
$ws.Range("D31").Value = '121.18'
$ws.Range("E31").Value = '  -1.19%  '

# Row 32
$ws.Range("D32").Value = '0.09509'
$ws.Range("E32").Value = '  -0.51%  '

# Row 33
$ws.Range("D33").Value = '0.9599'
$ws.Range("E33").Value = '  -1.43%  '

# Row 34
$ws.Range("D34").Value = '3.664'
$ws.Range("E34").Value = '  +1.07%  '

# Row 35
$ws.Range("D35").Value = '5.312'
$ws.Range("E35").Value = '  -0.10%  '

# Row 36
$ws.Range("D36").Value = '1.345'
$ws.Range("E36").Value = '  -7.49%  '

# Row 37
$ws.Range("D37").Value = '0.06132'
$ws.Range("E37").Value = '  +0.05%  '

# Row 38
$ws.Range("D38").Value = '0.02237'
$ws.Range("E38").Value = '  -1.08%  '

# Row 39
$ws.Range("D39").Value = '1.219'
$ws.Range("E39").Value = '  -0.94%  '

# Row 40
$ws.Range("B40").Value = 'Frax'
$ws.Range("C40").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D40").Value = '1.016'
$ws.Range("E40").Value = '  +1.47%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '8.130'
$ws.Range("E41").Value = '  -0.22%  '

# Row 42
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.5908'
$ws.Range("E42").Value = '  -1.75%  '

# Row 43
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '0.1887'
$ws.Range("E43").Value = '  -0.87%  '

# Row 44
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '10.15'
$ws.Range("E44").Value = '  -0.79%  '

# Row 45
$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '1.270'
$ws.Range("E45").Value = '  +2.04%  '

# Row 46
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.5638'
$ws.Range("E46").Value = '  -1.36%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '12.20'
$ws.Range("E47").Value = '  -0.69%  '

# Row 48
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").Value = '3.390'
$ws.Range("E48").Value = '  -0.29%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.917'
$ws.Range("E49").Value = '  -0.89%  '

# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.06885'
$ws.Range("E50").Value = '  +0.99%  '

# Row 51
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = '113.62'
$ws.Range("E51").Value = '  +1.71%  '

# Restore the default (no explicit number-format) style on the price column
# so the cells match the original workbook's formatting (no style override).
$priceRange.Style = "Normal"
